# Append a new paragraph "Change in other branch" after the existing
# "Second line" paragraph, and move the "_GoBack" bookmark (currently
# sitting right after "Second line") onto the end of this new paragraph.

$d = $word.ActiveDocument

# Locate the existing "_GoBack" bookmark; the new paragraph must be
# inserted exactly at its position so it lands right after "Second line"
# and before whatever used to follow the bookmark.
$goBack = $d.Bookmarks("_GoBack")
$insertAt = $goBack.Start

$insertionPoint = $d.Range($insertAt, $insertAt)
$insertionPoint.InsertParagraphAfter()

# The freshly inserted (still empty) paragraph is now the last one.
$newPara = $d.Paragraphs.Last

# Write the paragraph text with a one-character placeholder tacked on
# the end. Word COM's Bookmarks.Add has trouble placing a zero-length
# bookmark exactly at a paragraph's trailing mark (Range.End - 1); by
# temporarily keeping one extra character after the target spot, the
# bookmark is anchored at a safe, non-boundary offset first.
$newPara.Range.Text = "Change in other branchX"

$newParaRange = $d.Paragraphs.Last.Range
$bookmarkPos = $newParaRange.End - 2

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the placeholder character; it sits right after the bookmark so
# deleting it doesn't disturb the bookmark's now-settled position.
$placeholder = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholder.Delete()
